$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "86.272.72"
Set-TextValue $ws.Range("E2") "  +5.44%  "
Set-TextValue $ws.Range("D3") "3.267.76"
Set-TextValue $ws.Range("E3") "  +3.40%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "209.93"
Set-TextValue $ws.Range("E5") "  -2.98%  "
Set-TextValue $ws.Range("D6") "625.11"
Set-TextValue $ws.Range("E6") "  +1.36%  "
Set-TextValue $ws.Range("D7") "0.374"
Set-TextValue $ws.Range("E7") "  +29.81%  "
Set-TextValue $ws.Range("B8") "XRP"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D8") "0.651"
Set-TextValue $ws.Range("E8") "  +12.10%  "
Set-TextValue $ws.Range("B9") "USDC"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("E9") "  +0.02%  "
Set-TextValue $ws.Range("D10") "3.267.86"
Set-TextValue $ws.Range("E10") "  +3.52%  "
Set-TextValue $ws.Range("D11") "0.579"
Set-TextValue $ws.Range("E11") "  -3.79%  "
Set-TextValue $ws.Range("D12") "0.178"
Set-TextValue $ws.Range("E12") "  +8.16%  "
Set-TextValue $ws.Range("E13") "  +0.76%  "
Set-TextValue $ws.Range("B14") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D14") "3.859.00"
Set-TextValue $ws.Range("E14") "  +3.15%  "
Set-TextValue $ws.Range("B15") "Avalanche"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D15") "33.82"
Set-TextValue $ws.Range("E15") "  +6.06%  "
Set-TextValue $ws.Range("D16") "5.26"
Set-TextValue $ws.Range("E16") "  -0.52%  "
Set-TextValue $ws.Range("D17") "86.112.13"
Set-TextValue $ws.Range("E17") "  +5.44%  "
Set-TextValue $ws.Range("D18") "3.250.69"
Set-TextValue $ws.Range("E18") "  +2.91%  "
Set-TextValue $ws.Range("D19") "14.04"
Set-TextValue $ws.Range("E19") "  +0.48%  "
Set-TextValue $ws.Range("D20") "2.98"
Set-TextValue $ws.Range("E20") "  -7.17%  "
Set-TextValue $ws.Range("D21") "429.57"
Set-TextValue $ws.Range("E21") "  -0.90%  "
Set-TextValue $ws.Range("D22") "8.98"
Set-TextValue $ws.Range("E22") "  +1.13%  "
Set-TextValue $ws.Range("D23") "5.31"
Set-TextValue $ws.Range("E23") "  +4.16%  "
Set-TextValue $ws.Range("E24") "  -1.48%  "
Set-TextValue $ws.Range("D25") "12.25"
Set-TextValue $ws.Range("E25") "  +4.36%  "
Set-TextValue $ws.Range("D26") "5.14"
Set-TextValue $ws.Range("E26") "  -1.94%  "
Set-TextValue $ws.Range("D27") "3.426.96"
Set-TextValue $ws.Range("E27") "  +2.80%  "
Set-TextValue $ws.Range("D28") "75.86"
Set-TextValue $ws.Range("E28") "  -0.73%  "
Set-TextValue $ws.Range("D29") "0.0000129"
Set-TextValue $ws.Range("E29") "  +6.51%  "
Set-TextValue $ws.Range("E30") "  -0.03%  "
Set-TextValue $ws.Range("D31") "0.172"
Set-TextValue $ws.Range("E31") "  +17.38%  "
Set-TextValue $ws.Range("D32") "0.997"
Set-TextValue $ws.Range("E32") "  -0.36%  "
Set-TextValue $ws.Range("E33") "  -1.58%  "
Set-TextValue $ws.Range("D34") "544.17"
Set-TextValue $ws.Range("E34") "  -3.96%  "
Set-TextValue $ws.Range("E35") "  -4.86%  "
Set-TextValue $ws.Range("E36") "  -1.30%  "
Set-TextValue $ws.Range("D37") "6.82"
Set-TextValue $ws.Range("E37") "  +12.15%  "
Set-TextValue $ws.Range("E38") "  -10.12%  "
Set-TextValue $ws.Range("E39") "  -0.61%  "
Set-TextValue $ws.Range("D40") "0.997"
Set-TextValue $ws.Range("E40") "  -0.16%  "
Set-TextValue $ws.Range("D41") "21.64"
Set-TextValue $ws.Range("E41") "  +3.86%  "
Set-TextValue $ws.Range("D42") "0.393"
Set-TextValue $ws.Range("E42") "  -3.06%  "
Set-TextValue $ws.Range("E43") "  -0.56%  "
Set-TextValue $ws.Range("B44") "USDe"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D44") "1.00"
Set-TextValue $ws.Range("E44") "  -0.05%  "
Set-TextValue $ws.Range("B45") "dogwifhat"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D45") "2.91"
Set-TextValue $ws.Range("E45") "  -3.65%  "
Set-TextValue $ws.Range("D46") "155.70"
Set-TextValue $ws.Range("E46") "  -1.91%  "
Set-TextValue $ws.Range("D47") "177.07"
Set-TextValue $ws.Range("E47") "  -4.78%  "
Set-TextValue $ws.Range("E48") "  +0.16%  "
Set-TextValue $ws.Range("D49") "43.94"
Set-TextValue $ws.Range("E49") "  -1.77%  "
Set-TextValue $ws.Range("D50") "4.26"
Set-TextValue $ws.Range("E50") "  +2.09%  "
Set-TextValue $ws.Range("B51") "Mantle"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D51") "0.737"
Set-TextValue $ws.Range("E51") "  -3.13%  "
